# Updates the LR-pairs TPM results table on Sheet1 with newly computed
# values ("update scripts wuth new tpm"). The sending/target cluster set
# gained "ECs" (now 3 clusters: ECs, FAPs, MuSCs), so the data block grows
# from the original 3x2 (6 rows) combinations to a full 3x3 (9 rows) cross
# product of Sending cluster x Target cluster for the Clec11a -> Itga10
# ligand/receptor pair, occupying rows 2-10 (dimension becomes A1:T10).
# Row 1 (column headers) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Clec11a"
$ws.Cells.Item(2, 3).Value = "Itga10"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1.0
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.022591
$ws.Cells.Item(2, 8).Value = 0.067773
$ws.Cells.Item(2, 9).Value = 0.001469689085715816
$ws.Cells.Item(2, 10).Value = 0.001469689085715816
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 0.240998
$ws.Cells.Item(2, 14).Value = 0.7229939999999999
$ws.Cells.Item(2, 15).Value = 0.05495977716704094
$ws.Cells.Item(2, 16).Value = 0.05495977716704094
$ws.Cells.Item(2, 17).Value = 0.005444385817999999
$ws.Cells.Item(2, 18).Value = 0.04899947236199999
$ws.Cells.Item(2, 19).Value = 0.00008077378465577337
$ws.Cells.Item(2, 20).Value = 0.00008077378465577337

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Clec11a"
$ws.Cells.Item(3, 3).Value = "Itga10"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1.0
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.022591
$ws.Cells.Item(3, 8).Value = 0.067773
$ws.Cells.Item(3, 9).Value = 0.001469689085715816
$ws.Cells.Item(3, 10).Value = 0.001469689085715816
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 1.557543666666667
$ws.Cells.Item(3, 14).Value = 4.672631
$ws.Cells.Item(3, 15).Value = 0.3551990176181375
$ws.Cells.Item(3, 16).Value = 0.3551990176181375
$ws.Cells.Item(3, 17).Value = 0.03518646897366667
$ws.Cells.Item(3, 18).Value = 0.316678220763
$ws.Cells.Item(3, 19).Value = 0.0005220321194503565
$ws.Cells.Item(3, 20).Value = 0.0005220321194503565

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Clec11a"
$ws.Cells.Item(4, 3).Value = "Itga10"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1.0
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.022591
$ws.Cells.Item(4, 8).Value = 0.067773
$ws.Cells.Item(4, 9).Value = 0.001469689085715816
$ws.Cells.Item(4, 10).Value = 0.001469689085715816
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 2.586447
$ws.Cells.Item(4, 14).Value = 7.759341
$ws.Cells.Item(4, 15).Value = 0.5898412052148215
$ws.Cells.Item(4, 16).Value = 0.5898412052148215
$ws.Cells.Item(4, 17).Value = 0.05843042417700001
$ws.Cells.Item(4, 18).Value = 0.525873817593
$ws.Cells.Item(4, 19).Value = 0.000866883181609686
$ws.Cells.Item(4, 20).Value = 0.000866883181609686

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Clec11a"
$ws.Cells.Item(5, 3).Value = "Itga10"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 14.768619
$ws.Cells.Item(5, 8).Value = 44.305857
$ws.Cells.Item(5, 9).Value = 0.9607931545923257
$ws.Cells.Item(5, 10).Value = 0.9607931545923256
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 0.240998
$ws.Cells.Item(5, 14).Value = 0.7229939999999999
$ws.Cells.Item(5, 15).Value = 0.05495977716704094
$ws.Cells.Item(5, 16).Value = 0.05495977716704094
$ws.Cells.Item(5, 17).Value = 3.559207641762
$ws.Cells.Item(5, 18).Value = 32.032868775858
$ws.Cells.Item(5, 19).Value = 0.05280497768001254
$ws.Cells.Item(5, 20).Value = 0.05280497768001254

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Clec11a"
$ws.Cells.Item(6, 3).Value = "Itga10"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 14.768619
$ws.Cells.Item(6, 8).Value = 44.305857
$ws.Cells.Item(6, 9).Value = 0.9607931545923257
$ws.Cells.Item(6, 10).Value = 0.9607931545923256
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 1.557543666666667
$ws.Cells.Item(6, 14).Value = 4.672631
$ws.Cells.Item(6, 15).Value = 0.3551990176181375
$ws.Cells.Item(6, 16).Value = 0.3551990176181375
$ws.Cells.Item(6, 17).Value = 23.002768988863
$ws.Cells.Item(6, 18).Value = 207.024920899767
$ws.Cells.Item(6, 19).Value = 0.3412727846454254
$ws.Cells.Item(6, 20).Value = 0.3412727846454254

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Clec11a"
$ws.Cells.Item(7, 3).Value = "Itga10"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 14.768619
$ws.Cells.Item(7, 8).Value = 44.305857
$ws.Cells.Item(7, 9).Value = 0.9607931545923257
$ws.Cells.Item(7, 10).Value = 0.9607931545923256
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 2.586447
$ws.Cells.Item(7, 14).Value = 7.759341
$ws.Cells.Item(7, 15).Value = 0.5898412052148215
$ws.Cells.Item(7, 16).Value = 0.5898412052148215
$ws.Cells.Item(7, 17).Value = 38.19825030669301
$ws.Cells.Item(7, 18).Value = 343.7842527602371
$ws.Cells.Item(7, 19).Value = 0.5667153922668877
$ws.Cells.Item(7, 20).Value = 0.5667153922668877

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Clec11a"
$ws.Cells.Item(8, 3).Value = "Itga10"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 0.5800683333333333
$ws.Cells.Item(8, 8).Value = 1.740205
$ws.Cells.Item(8, 9).Value = 0.03773715632195847
$ws.Cells.Item(8, 10).Value = 0.03773715632195847
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 0.240998
$ws.Cells.Item(8, 14).Value = 0.7229939999999999
$ws.Cells.Item(8, 15).Value = 0.05495977716704094
$ws.Cells.Item(8, 16).Value = 0.05495977716704094
$ws.Cells.Item(8, 17).Value = 0.1397953081966666
$ws.Cells.Item(8, 18).Value = 1.25815777377
$ws.Cells.Item(8, 19).Value = 0.002074025702372628
$ws.Cells.Item(8, 20).Value = 0.002074025702372628

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Clec11a"
$ws.Cells.Item(9, 3).Value = "Itga10"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 0.5800683333333333
$ws.Cells.Item(9, 8).Value = 1.740205
$ws.Cells.Item(9, 9).Value = 0.03773715632195847
$ws.Cells.Item(9, 10).Value = 0.03773715632195847
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 1.557543666666667
$ws.Cells.Item(9, 14).Value = 4.672631
$ws.Cells.Item(9, 15).Value = 0.3551990176181375
$ws.Cells.Item(9, 16).Value = 0.3551990176181375
$ws.Cells.Item(9, 17).Value = 0.9034817588172221
$ws.Cells.Item(9, 18).Value = 8.131335829355
$ws.Cells.Item(9, 19).Value = 0.01340420085326173
$ws.Cells.Item(9, 20).Value = 0.01340420085326173

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Clec11a"
$ws.Cells.Item(10, 3).Value = "Itga10"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 0.5800683333333333
$ws.Cells.Item(10, 8).Value = 1.740205
$ws.Cells.Item(10, 9).Value = 0.03773715632195847
$ws.Cells.Item(10, 10).Value = 0.03773715632195847
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 2.586447
$ws.Cells.Item(10, 14).Value = 7.759341
$ws.Cells.Item(10, 15).Value = 0.5898412052148215
$ws.Cells.Item(10, 16).Value = 0.5898412052148215
$ws.Cells.Item(10, 17).Value = 1.500316000545
$ws.Cells.Item(10, 18).Value = 13.502844004905
$ws.Cells.Item(10, 19).Value = 0.02225892976632411
$ws.Cells.Item(10, 20).Value = 0.02225892976632411
